$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) cells we touch keep their original plain-text
# representation (e.g. '277.90', '1.003', '20.432.72') instead of being
# auto-converted to numbers by Excel, which would drop trailing zeros / grouping dots.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.432.72"
$ws.Range("E2").Value = "  -7.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.443.45"
$ws.Range("E3").Value = "  -6.98%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.003"
$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "277.90"
$ws.Range("E6").Value = "  -3.80%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3752"
$ws.Range("E7").Value = "  -4.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3069"
$ws.Range("E8").Value = "  -3.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "40.54"
$ws.Range("E9").Value = "  -8.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.014"
$ws.Range("E10").Value = "  -4.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06566"
$ws.Range("E11").Value = "  -7.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.16%  "

$ws.Range("E13").Value = "  -4.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.31"
$ws.Range("E14").Value = "  -6.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.140"
$ws.Range("E15").Value = "  -6.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.442.94"
$ws.Range("E16").Value = "  -7.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001011"
$ws.Range("E17").Value = "  -6.95%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05883"
$ws.Range("E18").Value = "  -10.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "76.33"
$ws.Range("E19").Value = "  -7.70%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.733"
$ws.Range("E21").Value = "  -6.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.42"
$ws.Range("E22").Value = "  -5.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.90"
$ws.Range("E23").Value = "  -1.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.311"
$ws.Range("E24").Value = "  -2.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.433.76"
$ws.Range("E25").Value = "  -7.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "143.18"
$ws.Range("E26").Value = "  -2.62%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.228"
$ws.Range("E27").Value = "  -5.37%  "

$ws.Range("E28").Value = "  -7.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.606.82"
$ws.Range("E29").Value = "  -6.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "109.54"
$ws.Range("E30").Value = "  -6.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.825"
$ws.Range("E31").Value = "  -21.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9085"
$ws.Range("E32").Value = "  -6.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.435"
$ws.Range("E33").Value = "  -6.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07740"
$ws.Range("E34").Value = "  -6.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.311"
$ws.Range("E35").Value = "  -8.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.003"
$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.84"
$ws.Range("E37").Value = "  +2.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05666"
$ws.Range("E38").Value = "  -4.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.145"
$ws.Range("E39").Value = "  -4.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.736"
$ws.Range("E40").Value = "  -6.11%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02047"
$ws.Range("E41").Value = "  -8.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1917"
$ws.Range("E42").Value = "  -5.18%  "

$ws.Range("B43").Value = "WEMIXTOKEN"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.359"
$ws.Range("E43").Value = "  -15.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.593"
$ws.Range("E44").Value = "  -3.90%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5336"
$ws.Range("E45").Value = "  -6.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.14"
$ws.Range("E46").Value = "  -5.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5169"
$ws.Range("E47").Value = "  -5.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "111.85"
$ws.Range("E48").Value = "  -3.89%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.791"
$ws.Range("E49").Value = "  -3.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.056"
$ws.Range("E50").Value = "  -5.85%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06291"
$ws.Range("E51").Value = "  -7.50%  "

